# Updated cryptos list with latest Price / Volume(1h) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.569.25"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "2.045.59"
$ws.Range("E3").Value = "  +3.36%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "256.53"
$ws.Range("E5").Value = "  +4.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.46"
$ws.Range("E8").Value = "  -5.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.14"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.82"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").Value = "2.346.77"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.821"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.43"
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.40"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "2.045.48"
$ws.Range("E18").Value = "  +3.27%  "
$ws.Range("D19").Value = "37.514.00"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.15"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.37"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  +7.48%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.21"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  -4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.54"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.98"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.35"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.77"
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0668"
$ws.Range("E34").Value = "  +7.44%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.50"
$ws.Range("E36").Value = "  +9.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.47"
$ws.Range("E37").Value = "  +4.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("E41").Value = "  +4.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0969"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.29"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "1.410.22"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.59"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.43"
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.04"
$ws.Range("E51").Value = "  +5.38%  "
